# Generate Report for Archive
#
# The localization-status report is regenerated: the "Ready for handoff"
# status text becomes "In Translation" everywhere it is used (the Overview
# sheet's per-locale status columns E/F, and the "Status" column on each
# per-locale detail sheet), and the now-narrower status columns are resized
# to fit the new (shorter) text.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"
$newColumnWidth = 12.5

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rowCount = $used.Rows.Count()
    $colCount = $used.Columns.Count()

    for ($r = 1; $r -le $rowCount; $r++) {
        for ($c = 1; $c -le $colCount; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            $value = $cell.Value()

            if ($oldStatus -eq $value) {
                $cell.Value = $newStatus
                $ws.Columns.Item($c).ColumnWidth = $newColumnWidth
            }
        }
    }
}
